$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (shared string text changes + date)
$ws.Range("A2").Value = "CDB"
$ws.Range("B2").Value = "CDB424EERW3"
$ws.Range("C2").Value = 46146

# Add new row 3 (introduce "CDB3239C4TL" before "BBRASIL FIM" so shared
# string indices line up: 6=CDB, 7=CDB424EERW3, 8=CDB3239C4TL, 9=BBRASIL FIM)
$ws.Range("A3").Value = "CDB"
$ws.Range("B3").Value = "CDB3239C4TL"
$ws.Range("C3").Value = 46223
$ws.Range("C3").NumberFormat = $ws.Range("C2").NumberFormat
$ws.Range("D3").Value = "BBRASIL FIM"

# Now update D2 to reuse the "BBRASIL FIM" shared string
$ws.Range("D2").Value = "BBRASIL FIM"
